# Update "countries & provincias Spain" data
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 26 de Mayo de 2020 a las 20:35"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1714371
$ws.Range("C4").Value = 8145
$ws.Range("E4").Value = 1145219

# Row 10 - Francia
$ws.Range("D10").Value = 65879
$ws.Range("E10").Value = 88533
$ws.Range("G10").Value = 98
$ws.Range("H10").Value = 28530

# Row 11 - Alemania
$ws.Range("B11").Value = 181203
$ws.Range("C11").Value = 414
$ws.Range("E11").Value = 10733
$ws.Range("G11").Value = 42
$ws.Range("H11").Value = 8470

# Row 16 - Canada
$ws.Range("B16").Value = 86614
$ws.Range("C16").Value = 903
$ws.Range("D16").Value = 45245
$ws.Range("E16").Value = 34732
$ws.Range("G16").Value = 92
$ws.Range("H16").Value = 6637

# Row 55 - Kazajistan
$ws.Range("D55").Value = 4613
$ws.Range("E55").Value = 4319
$ws.Range("G55").Value = 2
$ws.Range("H55").Value = 37

# Row 61 - Marruecos
$ws.Range("B61").Value = 7577
$ws.Range("C61").Value = 45
$ws.Range("D61").Value = 4881
$ws.Range("E61").Value = 2494

# Row 76 - Uzbekistan
$ws.Range("B76").Value = 3290
$ws.Range("C76").Value = 101
$ws.Range("E76").Value = 640
$ws.Range("G76").Value = 1
$ws.Range("H76").Value = 14

# Row 103 - Sri Lanka
$ws.Range("B103").Value = 1319
$ws.Range("C103").Value = 137
$ws.Range("E103").Value = 597

# Row 143 - Ruanda
$ws.Range("B143").Value = 339
$ws.Range("C143").Value = 3
$ws.Range("D143").Value = 244
$ws.Range("E143").Value = 95
